# Add a new "Potion" entity to the data dictionary / dependency matrix.
#
# Sheet "Dico de donnée" (sheet1) gets 4 new dictionary rows (ID_Potion,
# Nom_Potion, Prix_Potion, Temperature).
# Sheet "Feuil2" (sheet2) gets the matching 4 new rows/columns in the
# functional-dependency matrix, with a filled diagonal cell per new field
# and an "x" mark where Nom_Potion / Prix_Potion / Temperature depend on
# ID_Potion.

$wb = $excel.ActiveWorkbook
$wsDico = $wb.Worksheets.Item("Dico de donnée")
$wsMatrix = $wb.Worksheets.Item("Feuil2")

# ---------------------------------------------------------------
# Sheet "Dico de donnée": append the 4 new field definitions.
# ---------------------------------------------------------------
$dicoRows = @(
    @(45, "ID_Potion",   "Numerique",      "Obligatoire"),
    @(46, "Nom_Potion",  "Alphanumérique", "Obligatoire"),
    @(47, "Prix_Potion", "Monetaire",      "Obligatoire"),
    @(48, "Temperature", "Numerique",      "Obligatoire")
)

$r = 46
foreach ($row in $dicoRows) {
    $wsDico.Cells.Item($r, 1).Value = $row[0]
    $wsDico.Cells.Item($r, 2).Value = $row[1]
    $wsDico.Cells.Item($r, 3).Value = $row[2]
    $wsDico.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Feuil2": extend the header row with the 4 new column numbers.
# ---------------------------------------------------------------
$wsMatrix.Cells.Item(1, 47).Value = 45
$wsMatrix.Cells.Item(1, 48).Value = 46
$wsMatrix.Cells.Item(1, 49).Value = 47
$wsMatrix.Cells.Item(1, 50).Value = 48

# Row labels (same names/order as the dictionary sheet).
$wsMatrix.Cells.Item(46, 1).Value = 45
$wsMatrix.Cells.Item(46, 2).Value = "ID_Potion"
$wsMatrix.Cells.Item(47, 1).Value = 46
$wsMatrix.Cells.Item(47, 2).Value = "Nom_Potion"
$wsMatrix.Cells.Item(48, 1).Value = 47
$wsMatrix.Cells.Item(48, 2).Value = "Prix_Potion"
$wsMatrix.Cells.Item(49, 1).Value = 48
$wsMatrix.Cells.Item(49, 2).Value = "Temperature"

# Dependency marks: Nom_Potion / Prix_Potion / Temperature depend on ID_Potion
# (column AU, the diagonal column of row 46).
$wsMatrix.Cells.Item(47, 47).Value = "x"
$wsMatrix.Cells.Item(48, 47).Value = "x"
$wsMatrix.Cells.Item(49, 47).Value = "x"

# Diagonal "self" cells, filled black like the rest of the matrix
# (copy the format from an existing diagonal cell, e.g. C2).
$fmtSource = $wsMatrix.Range("C2")

$fmtSource.Copy()
$wsMatrix.Range("AU46").PasteSpecial(-4122)

$fmtSource.Copy()
$wsMatrix.Range("AV47").PasteSpecial(-4122)

$fmtSource.Copy()
$wsMatrix.Range("AW48").PasteSpecial(-4122)

$fmtSource.Copy()
$wsMatrix.Range("AX49").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Final view state: "Feuil2" keeps its own selection on the last new
# diagonal cell, but "Dico de donnée" becomes the active tab with the
# newly added rows selected.
# ---------------------------------------------------------------
$wsMatrix.Activate()
$wsMatrix.Range("AU49").Select()

$wsDico.Activate()
$wsDico.Range("B46:B49").Select()
